# Manual notes till feb 6th
# - Fill in the two new Adactin test scenario rows on the original sheet
# - Rename the original sheet to "Test Scenarios for Adactin"
# - Add a second sheet "Test Scenarios for Facebook" with its own scenarios
# - Re-point selections / column widths to match the final layout

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Fill in the two new rows on the Adactin (original) sheet ---
$ws1.Range("A7").Value = "Adactin_LoginPage_ForgotPassword_TS001"
$ws1.Range("B7").Value = "Login"
$ws1.Range("C7").Value = "R1.2"
$ws1.Range("D7").Value = "Forgot password functionailty"
$ws1.Range("E7").Value = "High"

$ws1.Range("A8").Value = "Adactin_LoginPage_Registration _TS001"
$ws1.Range("B8").Value = "Registration"
$ws1.Range("C8").Value = "R1.3"
$ws1.Range("D8").Value = "Registration functionailty"
$ws1.Range("E8").Value = "Very High"

# Rename the first sheet
$ws1.Name = "Test Scenarios for Adactin"

# --- Create the Facebook sheet after the Adactin sheet ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Test Scenarios for Facebook"

# Bring over the document-header rows (Project Name / References / Prepared / Reviewed / column titles)
$ws1.Range("A1:E5").Copy($ws2.Range("A1"))
# Stamp a clean bordered template row across the body so every data/blank row gets
# the same formatting (borders/fill) the Adactin sheet uses for its body rows.
$ws1.Range("A10:E10").Copy($ws2.Range("A6:E18"))

$facebookData = @(
  @("Facebook_LoginPage_Login_TS001","Login","R1.1","Login functionality using email","Very High"),
  @("Facebook_LoginPage_Login_TS002","Login","R1.2","Login functionality using mobile","Very High"),
  @("Facebook_LoginPage_ForgotPassword_TS001","Login","R1.3","Forgot password using email","Very High"),
  @("Facebook_LoginPage_ForgotPassword_TS002","Login","R1.4","Forgot password using mobile","Very High"),
  @("Facebook_LoginPage_CreateNewAccount_TS001","Registration","R1.5","Create a New Account Funtionality using Email","Very High"),
  @("Facebook_LoginPage_CreateNewAccount_TS002","Registration","R1.6","Create a New Account Funtionality using Mobile","Very High"),
  @("Facebook_LoginPage_CreatePage_TS001","Registration","R1.7","Create a Page Funtionality for Business or brand  ","Very High"),
  @("Facebook_LoginPage_CreatePage_TS002","Registration","R1.8","Create a Page Funtionality for Community or public figure ","Very High")
)

$r = 6
foreach ($row in $facebookData) {
  $ws2.Range("A$r").Value = $row[0]
  $ws2.Range("B$r").Value = $row[1]
  $ws2.Range("C$r").Value = $row[2]
  $ws2.Range("D$r").Value = $row[3]
  $ws2.Range("E$r").Value = $row[4]
  $r++
}

# --- Column widths ---
$ws1.Columns.Item(1).ColumnWidth = 54.45
$ws1.Columns.Item(4).ColumnWidth = 37.15

$ws2.Columns.Item(1).ColumnWidth = 60.45
$ws2.Columns.Item(2).ColumnWidth = 17.8
$ws2.Columns.Item(3).ColumnWidth = 18.6
$ws2.Columns.Item(4).ColumnWidth = 71.95
$ws2.Columns.Item(5).ColumnWidth = 12.0

# --- Row heights on the new sheet (match the 21pt rows used on the Adactin sheet) ---
for ($i = 1; $i -le 18; $i++) {
  $ws2.Rows.Item($i).RowHeight = 21
}

# --- Selections / active sheet ---
$ws1.Activate()
$ws1.Range("A6:E18").Select()

$ws2.Activate()
$ws2.Range("D15").Select()
